# Append a new data row (row 36) to each of the 4 worksheets, mirroring
# the layout/format of the existing last row (row 35), and grow the
# used range accordingly (A1:I35 -> A1:I36 on every sheet).

$wb = $excel.ActiveWorkbook

# time (column A) is identical for the new row on every sheet.
$newTime = 45822.49130787037

# Per-sheet values for the new row 36: B..I
$sheetData = @{
    1 = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x6C"
        E = "0xf"
        F = 380
        G = 759863127514710945038336.0
        H = 364
        I = 15
    }
    2 = @{
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x80"
        E = "0xe"
        F = 400
        G = 568432987514711010443264.0
        H = 384
        I = 14
    }
    3 = @{
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x6B"
        E = "0x3"
        F = 110
        G = 568631262647113970876416.0
        H = 107
        I = 3
    }
    4 = @{
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x6C"
        E = "0x3"
        F = 110
        G = 985046333984776009023488.0
        H = 108
        I = 3
    }
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $row = $sheetData[$i]

    # New row is one below the previous last data row (35 -> 36).
    $ws.Cells.Item(36, 1).Value = $newTime
    # Match the date/time number format used by the existing column-A cells.
    $ws.Cells.Item(36, 1).NumberFormat = $ws.Cells.Item(35, 1).NumberFormat

    $ws.Cells.Item(36, 2).Value = $row.B
    $ws.Cells.Item(36, 3).Value = $row.C
    $ws.Cells.Item(36, 4).Value = $row.D
    $ws.Cells.Item(36, 5).Value = $row.E
    $ws.Cells.Item(36, 6).Value = $row.F
    $ws.Cells.Item(36, 7).Value = $row.G
    $ws.Cells.Item(36, 8).Value = $row.H
    $ws.Cells.Item(36, 9).Value = $row.I
}
